$wb = $excel.ActiveWorkbook

# Sheet that currently holds the rows we need to split.
$userSheet = $wb.Worksheets.Item("user")

# Insert a brand-new sheet right after "user" - this becomes the new
# "Sheet1" tab that receives the row that used to be row 2 of "user".
$newSheet = $wb.Worksheets.Add($null, $userSheet)

# Move (cut/paste) the second data row of "user" (A2:R2, Ritu's record)
# onto row 1 of the new sheet - this carries cell values AND formatting
# (styles s="4"/s="5") across intact.
$userSheet.Range("A2:R2").Cut($newSheet.Range("A1"))

# Remove the now-empty row 2 from "user" so the remaining data row
# (previously row 3) shifts up to become row 2.
$userSheet.Rows(2).Delete()

# The hyperlinks that lived on (old) row 3 need to be re-pointed at the
# cells' new address (row 2). Re-create them there.
$userSheet.Hyperlinks.Delete()
$userSheet.Hyperlinks.Add($userSheet.Range("E2"), "mailto:nevixo9520@ociun.com")
$userSheet.Hyperlinks.Add($userSheet.Range("D2"), "mailto:String@123")
# Restore the plain "Hyperlink" cell style (Add() above re-styled the
# cells with a fresh style entry) so the cells match their original look.
$userSheet.Range("D2").Style = "Hyperlink"
$userSheet.Range("E2").Style = "Hyperlink"

# Mark the whole first row as selected on the new sheet (matches a
# "select entire row" gesture after pasting the cut data there).
$newSheet.Rows(1).Select()

# Finally put the active selection back on "user", which becomes the
# active tab, with the cursor at B16.
$userSheet.Activate()
$userSheet.Range("B16").Select()
